# Adds a second (fake/test) case row to the "Case Info" and "CSV Case Info"
# sheets, as described in the commit message:
#   "Added a second fake case to csv tab and case tab"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Data for the new case (row 3 on "Case Info", row 2 on "CSV Case Info").
# Keys are column letters, values are plain strings/numbers as typed by a
# user filling out the form.
# ---------------------------------------------------------------------------

function Fill-CaseRow($ws, $r) {
    $ws.Range("A$r").Value = 2
    $ws.Range("B$r").Value = "Fak E. Case"
    $ws.Range("C$r").Value = "winky eyes"
    $ws.Range("D$r").Value = "Jimmy Simonse"

    # Dates (formatted like the existing rows, m/d)
    $ws.Range("E2").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)
    $ws.Range("E$r").Value = 43792
    $ws.Range("F2").Copy()
    $ws.Range("F$r").PasteSpecial(-4122)
    $ws.Range("F$r").Value = 43792

    $ws.Range("G$r").Value = "test2"
    $ws.Range("H$r").Value = "none"

    $ws.Range("I$r").Value = "jsimonse9@gmail.com"
    $ws.Hyperlinks.Add($ws.Range("I$r"), "mailto:jsimonse9@gmail.com")
    $ws.Range("I$r").Style = "Hyperlink"

    $ws.Range("L$r").Value = 937
    $ws.Range("M$r").Value = "female"
    $ws.Range("N$r").Value = "2 kg"

    $ws.Range("O$r").Value = "A dark and shady history 1"
    $ws.Range("P$r").Value = "A darker and shadier history 2"

    $ws.Range("Q$r").Value = "Cat cuddles"
    $ws.Range("R$r").Value = "Adderol"
    $ws.Range("S$r").Value = "Vodka"

    $ws.Range("T$r").Value = "Wears oversized pants"
    $ws.Range("U$r").Value = "No rhinoceruses`n"
    $ws.Range("U$r").WrapText = $true
    $ws.Range("U$r").VerticalAlignment = -4160

    $ws.Range("V$r").Value = "No cough, no SOB"
    $ws.Range("W$r").Value = "Ny algae"
    $ws.Range("X$r").Value = "Ready to diffuse"
    $ws.Range("Y$r").Value = "Vampire bites"
    $ws.Range("Z$r").Value = "Brusing on nose hairs"
    $ws.Range("AA$r").Value = "Constantly drinks vodka"
    $ws.Range("AB$r").Value = "Anxious about clowns"

    $ws.Range("AC$r").Value = "Paper like thinness"
    $ws.Range("AD$r").Value = "pettasatus"
    $ws.Range("AE$r").Value = "Supple nose and left ear"
    $ws.Range("AF$r").Value = "Normal S1, S2, S3"
    $ws.Range("AG$r").Value = "Has two lungs"
    $ws.Range("AH$r").Value = "Has two arms"

    $ws.Range("AI$r").Value = "Periumbilical TTP with guarding, + McBurney point tenderness, - Murphy's Sign "
    $ws.Range("AJ$r").Value = "No CVA tenderness BL, "
    $ws.Range("AK$r").Value = "normal skin color and turgor, no jaundice, pallor, rashes, bruising, erythema "
    $ws.Range("AL$r").Value = "Strength and Sensation +5/5 BL upper and lower extremities, sensation grossly intact"
    $ws.Range("AM$r").Value = "A&O x3, cooperative"

    $ws.Range("AN$r").Value = 37
    $ws.Range("AO$r").Value = 129
    $ws.Range("AP$r").Value = 69
    $ws.Range("AQ$r").Value = 99
    $ws.Range("AR$r").Value = 19
    $ws.Range("AS$r").Value = 98

    $ws.Range("AT$r").Value = "Dry Eyes"
    $ws.Range("AU$r").Value = "Labs1"
    $ws.Range("AV$r").Value = "Labs2"
    $ws.Range("AW$r").Value = "Labs3"
    $ws.Range("AX$r").Value = "Labs4"
    $ws.Range("AY$r").Value = "rule out "
}

# ---------------------------------------------------------------------------
# "Case Info" sheet: header in row 1, first (real) case in row 2, so the new
# fake case goes in row 3.
# ---------------------------------------------------------------------------
$wsCase = $wb.Worksheets.Item("Case Info")
Fill-CaseRow $wsCase 3

# ---------------------------------------------------------------------------
# "CSV Case Info" sheet: no header row, first case is row 1, so the new fake
# case goes in row 2.
# ---------------------------------------------------------------------------
$wsCsv = $wb.Worksheets.Item("CSV Case Info")
Fill-CaseRow $wsCsv 2

# ---------------------------------------------------------------------------
# Leave the workbook with the selection/active sheet similar to how it was
# left by the author after typing in the new data.
# ---------------------------------------------------------------------------
$wsCase.Select()
$wsCase.Range("A3:AY3").Select()

$wsCsv.Select()
$wsCsv.Range("D4").Select()
